$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.287.73"
$ws.Range("E2").Value = "  -3.97%  "

$ws.Range("D3").Value = "3.304.35"
$ws.Range("E3").Value = "  -4.25%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'560.09"
$ws.Range("E5").Value = "  -3.23%  "

$ws.Range("D6").Value = "'143.75"
$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "3.307.21"
$ws.Range("E8").Value = "  -4.15%  "

$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("E10").Value = "  -2.17%  "

$ws.Range("E11").Value = "  -3.93%  "

$ws.Range("D12").Value = "'0.407"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").Value = "3.875.53"
$ws.Range("E13").Value = "  -4.13%  "

$ws.Range("E14").Value = "  +0.51%  "

$ws.Range("D15").Value = "'27.31"
$ws.Range("E15").Value = "  -3.45%  "

$ws.Range("D16").Value = "3.300.03"
$ws.Range("E16").Value = "  -4.20%  "

$ws.Range("E17").Value = "  -4.07%  "

$ws.Range("D18").Value = "60.301.31"
$ws.Range("E18").Value = "  -3.99%  "

$ws.Range("D19").Value = "'6.12"
$ws.Range("E19").Value = "  -4.19%  "

$ws.Range("D20").Value = "'14.37"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").Value = "'8.60"
$ws.Range("E21").Value = "  -4.56%  "

$ws.Range("D22").Value = "'374.25"
$ws.Range("E22").Value = "  -3.41%  "

$ws.Range("D23").Value = "'73.73"
$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("D24").Value = "'0.549"
$ws.Range("E24").Value = "  -3.41%  "

$ws.Range("E25").Value = "  +0.17%  "

$ws.Range("D26").Value = "3.460.20"
$ws.Range("E26").Value = "  -3.51%  "

$ws.Range("E27").Value = "  -9.18%  "

$ws.Range("D28").Value = "'0.173"
$ws.Range("E28").Value = "  -5.42%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "'7.23"
$ws.Range("E30").Value = "  -6.35%  "

$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'7.63"
$ws.Range("E32").Value = "  -4.75%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'2.05"
$ws.Range("E33").Value = "  -3.37%  "

$ws.Range("D34").Value = "'22.58"
$ws.Range("E34").Value = "  -2.89%  "

$ws.Range("D35").Value = "'1.28"
$ws.Range("E35").Value = "  -4.53%  "

$ws.Range("D36").Value = "'5.19"
$ws.Range("E36").Value = "  -4.42%  "

$ws.Range("D37").Value = "'166.86"
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("E39").Value = "  -7.42%  "

$ws.Range("D40").Value = "'27.48"
$ws.Range("E40").Value = "  -14.36%  "

$ws.Range("D41").Value = "3.338.07"
$ws.Range("E41").Value = "  -4.22%  "

$ws.Range("E42").Value = "  -5.52%  "

$ws.Range("D43").Value = "'41.85"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("D44").Value = "'0.752"
$ws.Range("E44").Value = "  -4.32%  "

$ws.Range("E45").Value = "  -4.11%  "

$ws.Range("D46").Value = "'1.60"
$ws.Range("E46").Value = "  -5.41%  "

$ws.Range("D47").Value = "'1.13"
$ws.Range("E47").Value = "  -3.95%  "

$ws.Range("D48").Value = "2.370.54"
$ws.Range("E48").Value = "  -7.67%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").Value = "'6.58"
$ws.Range("E50").Value = "  -4.76%  "

$ws.Range("D51").Value = "'21.55"
$ws.Range("E51").Value = "  -4.85%  "
